# audit_2026_1.xlsx - add "Grade" column, refresh sample rows, append new row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at C ("Grade") - shifts old Supplier..Leakage% from C..I to D..J
$ws.Columns("C").Insert()

# 2. New header for the inserted column
$ws.Range("C1").Value = "Grade"

# Columns A (Date) and B (text that could look numeric) must stay text, so force
# a text number format before writing values that Excel would otherwise reinterpret
# (dates / pure-digit strings) as numbers.
$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"

# 3. Row 2
$ws.Range("A2").Value = "2026-01-14"
$ws.Range("B2").Value = "AS 01"
$ws.Range("C2").Value = "VG10"
$ws.Range("D2").Value = "IOCL"
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = 50000
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 250000
$ws.Range("J2").Value = 10

# 4. Row 3
$ws.Range("A3").Value = "2026-01-18"
$ws.Range("B3").Value = "TEST "
$ws.Range("C3").Value = "VG30"
$ws.Range("D3").Value = "IOCL"
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = 50000
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 50000
$ws.Range("J3").Value = 10

# 5. Row 4 (new)
$ws.Range("A4").Value = "2026-01-18"
$ws.Range("B4").Value = "123"
$ws.Range("C4").Value = "VG10"
$ws.Range("D4").Value = "IOML"
$ws.Range("E4").Value = 50
$ws.Range("F4").Value = 45
$ws.Range("G4").Value = 50000
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 250000
$ws.Range("J4").Value = 10

Write-Output "edit complete"
